$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "edit1"
$ws.Range("B10").Value = "riya-morankar"
$ws.Range("C10").Value = "Squashed"
$ws.Range("D10").Value = "2 changes"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2025-06-18"
$ws.Range("F10").Value = "N/A"
